# "요청 모니터링.xlsx" update — request monitoring sheet status refresh
# - Row 4 (request #682): status "요청" -> "진행중"
# - Row 5 (request #681): status "요청" -> "완료"; fill in start/end timestamps (H/I);
#   update the "요청일시" (L) timestamp
# - Row 6 (request #680): status "요청" -> "완료"; fill in start/end timestamps (H/I);
#   update the "요청일시" (L) timestamp
# - Row 7 (request #679): status "요청" -> "진행대기"; update the "요청일시" (L) timestamp
# - Move the active selection to K19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing date/time display format (column L already uses it) so the
# newly written H/I timestamp cells match the workbook's existing style.
$dateFormat = "m/d/yy h:mm"

# Row 4 — move to "진행중" (in progress)
$ws.Range("J4").Value = "진행중"

# Row 5 — move to "완료" (done); record start/end timestamps
$ws.Range("H5").Value = 45765.389884259261
$ws.Range("H5").NumberFormat = $dateFormat
$ws.Range("I5").Value = 45765.390879629631
$ws.Range("I5").NumberFormat = $dateFormat
$ws.Range("J5").Value = "완료"
$ws.Range("L5").Value = 45758.606458333335

# Row 6 — move to "완료" (done); record start/end timestamps
$ws.Range("H6").Value = 45764.860937500001
$ws.Range("H6").NumberFormat = $dateFormat
$ws.Range("I6").Value = 45764.863877314812
$ws.Range("I6").NumberFormat = $dateFormat
$ws.Range("J6").Value = "완료"
$ws.Range("L6").Value = 45758.696886574071

# Row 7 — move to "진행대기" (pending)
$ws.Range("J7").Value = "진행대기"
$ws.Range("L7").Value = 45758.744120370371

# Match the author's last on-screen selection
$ws.Range("K19").Select()
